$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D contains numeric-looking text values (prices) that must remain text.
# Apply a Text number format to the whole price column first so Excel does not
# auto-convert the assigned strings into numbers, then restore the default style.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '43.849.80'
$ws.Range("E2").Value = '  +0.17%  '

$ws.Range("D3").Value = '2.294.09'
$ws.Range("E3").Value = '  +0.30%  '

$ws.Range("E4").Value = '  +0.24%  '

$ws.Range("D5").Value = '113.41'
$ws.Range("E5").Value = '  +17.34%  '

$ws.Range("D6").Value = '269.49'
$ws.Range("E6").Value = '  -0.05%  '

$ws.Range("E7").Value = '  +1.62%  '

$ws.Range("E8").Value = '  +0.28%  '

$ws.Range("D9").Value = '0.617'
$ws.Range("E9").Value = '  +1.46%  '

$ws.Range("D10").Value = '48.01'
$ws.Range("E10").Value = '  +6.54%  '

$ws.Range("D11").Value = '0.0944'
$ws.Range("E11").Value = '  +1.44%  '

$ws.Range("E12").Value = '  +14.79%  '

$ws.Range("E13").Value = '  +0.90%  '

$ws.Range("D14").Value = '15.87'
$ws.Range("E14").Value = '  +0.28%  '

$ws.Range("D15").Value = '2.639.13'
$ws.Range("E15").Value = '  +0.32%  '

$ws.Range("D16").Value = '0.862'
$ws.Range("E16").Value = '  +0.51%  '

$ws.Range("D17").Value = '2.298.42'
$ws.Range("E17").Value = '  +0.35%  '

$ws.Range("D18").Value = '43.732.59'
$ws.Range("E18").Value = '  +0.00%  '

$ws.Range("E19").Value = '  -0.90%  '

$ws.Range("D20").Value = '6.82'
$ws.Range("E20").Value = '  +10.40%  '

$ws.Range("D21").Value = '72.20'
$ws.Range("E21").Value = '  +0.20%  '

$ws.Range("D22").Value = '2.43'
$ws.Range("E22").Value = '  -2.33%  '

$ws.Range("D23").Value = '3.02'
$ws.Range("E23").Value = '  +12.04%  '

$ws.Range("D24").Value = '233.03'
$ws.Range("E24").Value = '  +0.15%  '

$ws.Range("D25").Value = '9.66'
$ws.Range("E25").Value = '  +5.85%  '

$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  +0.03%  '

$ws.Range("D27").Value = '11.60'
$ws.Range("E27").Value = '  +2.17%  '

$ws.Range("D28").Value = '41.49'
$ws.Range("E28").Value = '  +7.33%  '

$ws.Range("D29").Value = '3.39'
$ws.Range("E29").Value = '  -1.92%  '

$ws.Range("E30").Value = '  -0.73%  '

$ws.Range("D31").Value = '175.56'
$ws.Range("E31").Value = '  +0.13%  '

$ws.Range("D32").Value = '0.0928'
$ws.Range("E32").Value = '  +3.87%  '

$ws.Range("D33").Value = '21.55'
$ws.Range("E33").Value = '  -1.14%  '

$ws.Range("E34").Value = '  +5.32%  '

$ws.Range("E35").Value = '  +0.16%  '

$ws.Range("E36").Value = '  -1.90%  '

$ws.Range("E37").Value = '  +3.57%  '

$ws.Range("E38").Value = '  +0.46%  '

$ws.Range("D39").Value = '3.82'
$ws.Range("E39").Value = '  +8.56%  '

$ws.Range("B40").Value = 'Celestia'
$ws.Range("C40").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D40").Value = '13.92'
$ws.Range("E40").Value = '  +13.53%  '

$ws.Range("B41").Value = 'MultiversX'
$ws.Range("C41").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D41").Value = '73.84'
$ws.Range("E41").Value = '  +14.43%  '

$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D42").Value = '0.243'
$ws.Range("E42").Value = '  +2.92%  '

$ws.Range("E43").Value = '  +3.06%  '

$ws.Range("D44").Value = '6.27'
$ws.Range("E44").Value = '  +20.70%  '

$ws.Range("E45").Value = '  +0.21%  '

$ws.Range("E46").Value = '  +3.42%  '

$ws.Range("D47").Value = '8.82'
$ws.Range("E47").Value = '  +1.37%  '

$ws.Range("E48").Value = '  -2.41%  '

$ws.Range("D49").Value = '101.45'
$ws.Range("E49").Value = '  +3.33%  '

$ws.Range("E50").Value = '  +3.23%  '

$ws.Range("D51").Value = '0.466'
$ws.Range("E51").Value = '  +5.58%  '

# Restore the default (Normal) style on column D so no extra formatting lingers
$ws.Range("D2:D51").Style = "Normal"
